$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain plain text (matching source inlineStr formatting)
# Force text number format first so Excel does not auto-convert them to numbers.
$textCells = @("D4","D5","D6","D11","D12","D13","D17","D19","D22","D24","D25","D27","D28","D31","D32","D33","D34","D35","D38","D39","D40","D41","D43","D46","D47","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values cell by cell, matching the source diff
$ws.Range("D2").Value = '67.348.70'
$ws.Range("D3").Value = '3.255.27'
$ws.Range("E3").Value = '  -7.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '593.12'
$ws.Range("E5").Value = '  -4.37%  '
$ws.Range("D6").Value = '150.85'
$ws.Range("E6").Value = '  -12.80%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.247.16'
$ws.Range("E8").Value = '  -7.42%  '
$ws.Range("E9").Value = '  -11.02%  '
$ws.Range("E10").Value = '  -13.34%  '
$ws.Range("D11").Value = '6.80'
$ws.Range("E11").Value = '  -3.38%  '
$ws.Range("D12").Value = '0.506'
$ws.Range("E12").Value = '  -12.96%  '
$ws.Range("D13").Value = '38.11'
$ws.Range("E13").Value = '  -17.63%  '
$ws.Range("E14").Value = '  -12.19%  '
$ws.Range("D15").Value = '3.771.92'
$ws.Range("E15").Value = '  -7.54%  '
$ws.Range("D16").Value = '67.333.34'
$ws.Range("E16").Value = '  -4.65%  '
$ws.Range("D17").Value = '546.20'
$ws.Range("E17").Value = '  -10.51%  '
$ws.Range("D18").Value = '3.256.16'
$ws.Range("E18").Value = '  -7.14%  '
$ws.Range("D19").Value = '7.24'
$ws.Range("E19").Value = '  -13.61%  '
$ws.Range("E20").Value = '  -5.99%  '
$ws.Range("E21").Value = '  -14.36%  '
$ws.Range("D22").Value = '0.762'
$ws.Range("E22").Value = '  -13.45%  '
$ws.Range("E23").Value = '  -14.05%  '
$ws.Range("D24").Value = '85.36'
$ws.Range("E24").Value = '  -13.15%  '
$ws.Range("D25").Value = '13.50'
$ws.Range("E25").Value = '  -13.10%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '3.23'
$ws.Range("E27").Value = '  -13.66%  '
$ws.Range("D28").Value = '8.04'
$ws.Range("E28").Value = '  -11.11%  '
$ws.Range("E29").Value = '  -13.03%  '
$ws.Range("E30").Value = '  -17.34%  '
$ws.Range("D31").Value = '2.65'
$ws.Range("E31").Value = '  -12.14%  '
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  -12.59%  '
$ws.Range("D33").Value = '542.61'
$ws.Range("E33").Value = '  -14.50%  '
$ws.Range("D34").Value = '6.63'
$ws.Range("E34").Value = '  -17.82%  '
$ws.Range("D35").Value = '5.70'
$ws.Range("E35").Value = '  -15.84%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  -5.17%  '
$ws.Range("D38").Value = '53.69'
$ws.Range("E38").Value = '  -5.43%  '
$ws.Range("D39").Value = '0.0854'
$ws.Range("E39").Value = '  -14.23%  '
$ws.Range("D40").Value = '9.16'
$ws.Range("E40").Value = '  -14.88%  '
$ws.Range("D41").Value = '0.126'
$ws.Range("E41").Value = '  -11.34%  '
$ws.Range("D42").Value = '2.930.33'
$ws.Range("E42").Value = '  -12.26%  '
$ws.Range("D43").Value = '2.62'
$ws.Range("E43").Value = '  -23.14%  '
$ws.Range("E44").Value = '  -16.17%  '
$ws.Range("D45").Value = '0.0₃0579'
$ws.Range("E45").Value = '  -19.14%  '
$ws.Range("D46").Value = '26.43'
$ws.Range("E46").Value = '  -16.75%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.14'
$ws.Range("E47").Value = '  -15.14%  '
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '2.34'
$ws.Range("E49").Value = '  -20.45%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '127.19'
$ws.Range("E50").Value = '  -4.94%  '
$ws.Range("E51").Value = '  -12.38%  '
